# no-op for now, just testing
$p = $ppt.ActivePresentation
